# Scheduled runner: refresh market-price/profit columns (H-N) across job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 87704.74000000001
$ws.Range("I28").Value = 95895.81
$ws.Range("J28").Value = 1698.5
$ws.Range("K28").Value = 95895.81
$ws.Range("L28").Value = 1698.5
$ws.Range("M28").Value = -95410.81
$ws.Range("N28").Value = -2668.5
$ws.Range("H40").Value = 3904.3044
$ws.Range("J40").Value = 3999.9524
$ws.Range("L40").Value = 3999.9524
$ws.Range("N40").Value = -4349.9524
$ws.Range("H74").Value = 9250.75
$ws.Range("I74").Value = 9250.75
$ws.Range("K74").Value = 9250.75
$ws.Range("M74").Value = -8314.75
$ws.Range("H77").Value = 9250.75
$ws.Range("I77").Value = 9250.75
$ws.Range("K77").Value = 46253.75
$ws.Range("M77").Value = -41573.75
$ws.Range("H86").Value = 250002320
$ws.Range("J86").Value = 250002080
$ws.Range("L86").Value = 250002080
$ws.Range("N86").Value = -250004326
$ws.Range("H89").Value = 250002320
$ws.Range("J89").Value = 250002080
$ws.Range("L89").Value = 1250010400
$ws.Range("N89").Value = -1250021632
$ws.Range("H103").Value = 1410.2858
$ws.Range("I103").Value = 1237.9474
$ws.Range("J103").Value = 1614.9375
$ws.Range("K103").Value = 3713.8422
$ws.Range("L103").Value = 4844.8125
$ws.Range("M103").Value = -3127.8422
$ws.Range("N103").Value = -6016.8125
$ws.Range("H112").Value = 879880.6
$ws.Range("J112").Value = 1152152.4
$ws.Range("L112").Value = 3456457.2
$ws.Range("N112").Value = -3458673.2
$ws.Range("H135").Value = 2402.7222
$ws.Range("I135").Value = 1095.9166
$ws.Range("J135").Value = 5016.3335
$ws.Range("K135").Value = 9863.249400000001
$ws.Range("L135").Value = 45147.0015
$ws.Range("M135").Value = -7328.249400000001
$ws.Range("N135").Value = -50217.0015
$ws.Range("H136").Value = 112628.57
$ws.Range("J136").Value = 112628.57
$ws.Range("L136").Value = 112628.57
$ws.Range("N136").Value = -122828.57
$ws.Range("H137").Value = 5558985.5
$ws.Range("I137").Value = 3333.1667
$ws.Range("K137").Value = 9999.500100000001
$ws.Range("M137").Value = -7449.500100000001
$ws.Range("H138").Value = 2513.6565
$ws.Range("I138").Value = 1220.1333
$ws.Range("J138").Value = 2744.6428
$ws.Range("K138").Value = 3660.3999
$ws.Range("L138").Value = 8233.928400000001
$ws.Range("M138").Value = 1479.6001
$ws.Range("N138").Value = -18513.9284
$ws.Range("H141").Value = 8021.8125
$ws.Range("I141").Value = 7178.0557
$ws.Range("K141").Value = 21534.1671
$ws.Range("M141").Value = -16354.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3576.2144
$ws.Range("I2").Value = 3408.5
$ws.Range("K2").Value = 3408.5
$ws.Range("M2").Value = -3295.5
$ws.Range("H21").Value = 27500
$ws.Range("I21").Value = 15000
$ws.Range("J21").Value = 40000
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 40000
$ws.Range("M21").Value = -14626
$ws.Range("N21").Value = -40748
$ws.Range("H32").Value = 27537022
$ws.Range("I32").Value = 29114174
$ws.Range("K32").Value = 29114174
$ws.Range("M32").Value = -29113887
$ws.Range("H74").Value = 3507.3125
$ws.Range("I74").Value = 3668.3845
$ws.Range("K74").Value = 3668.3845
$ws.Range("M74").Value = -2794.3845
$ws.Range("H77").Value = 3507.3125
$ws.Range("I77").Value = 3668.3845
$ws.Range("K77").Value = 18341.9225
$ws.Range("M77").Value = -13973.9225
$ws.Range("H110").Value = 1615.7142
$ws.Range("I110").Value = 373.33334
$ws.Range("K110").Value = 373.33334
$ws.Range("M110").Value = 1671.66666
$ws.Range("H116").Value = 3576.2144
$ws.Range("I116").Value = 3408.5
$ws.Range("K116").Value = 3408.5
$ws.Range("M116").Value = -1114.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3576.2144
$ws.Range("I3").Value = 3408.5
$ws.Range("K3").Value = 3408.5
$ws.Range("M3").Value = -3294.5
$ws.Range("H107").Value = 800.1905
$ws.Range("I107").Value = 601.6667
$ws.Range("J107").Value = 1991.3334
$ws.Range("K107").Value = 601.6667
$ws.Range("L107").Value = 1991.3334
$ws.Range("M107").Value = 1318.3333
$ws.Range("N107").Value = -5831.3334
$ws.Range("H124").Value = 179000
$ws.Range("J124").Value = 179000
$ws.Range("L124").Value = 179000
$ws.Range("N124").Value = -188820
$ws.Range("H126").Value = 136000
$ws.Range("J126").Value = 136000
$ws.Range("L126").Value = 136000
$ws.Range("N126").Value = -145880
$ws.Range("H134").Value = 2978702
$ws.Range("I134").Value = 3970371.8
$ws.Range("K134").Value = 11911115.4
$ws.Range("M134").Value = -11908580.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 500000
$ws.Range("I4").Value = 500000
$ws.Range("K4").Value = 500000
$ws.Range("M4").Value = -499888
$ws.Range("H125").Value = 91167
$ws.Range("J125").Value = 91167
$ws.Range("L125").Value = 91167
$ws.Range("N125").Value = -96087

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 755.3333
$ws.Range("I103").Value = 648.5
$ws.Range("K103").Value = 1945.5
$ws.Range("M103").Value = -1066.5
$ws.Range("H107").Value = 981.5
$ws.Range("I107").Value = 1108.7059
$ws.Range("K107").Value = 3326.1177
$ws.Range("M107").Value = -1406.1177
$ws.Range("H113").Value = 1346.3572
$ws.Range("J113").Value = 1489.2727
$ws.Range("L113").Value = 4467.8181
$ws.Range("N113").Value = -8807.8181
$ws.Range("H114").Value = 809.7778
$ws.Range("J114").Value = 698.5
$ws.Range("L114").Value = 2095.5
$ws.Range("N114").Value = -8603.5
$ws.Range("H121").Value = 5033958
$ws.Range("I121").Value = 780
$ws.Range("J121").Value = 5663105
$ws.Range("K121").Value = 2340
$ws.Range("L121").Value = 16989315
$ws.Range("M121").Value = -1030
$ws.Range("N121").Value = -16991935
$ws.Range("H122").Value = 107.833336
$ws.Range("J122").Value = 150
$ws.Range("L122").Value = 1350
$ws.Range("N122").Value = -6250
$ws.Range("H132").Value = 418330.22
$ws.Range("I132").Value = 1288.4615
$ws.Range("K132").Value = 11596.1535
$ws.Range("M132").Value = -9066.153499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 125464.8
$ws.Range("J134").Value = 125464.8
$ws.Range("L134").Value = 376394.4
$ws.Range("N134").Value = -381464.4
$ws.Range("H136").Value = 70712.2
$ws.Range("J136").Value = 70712.2
$ws.Range("L136").Value = 212136.6
$ws.Range("N136").Value = -217236.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1882501.5
$ws.Range("H82").Value = 2485.8572
$ws.Range("I82").Value = 2558.8667
$ws.Range("K82").Value = 2558.8667
$ws.Range("M82").Value = -2197.8667
$ws.Range("H85").Value = 2485.8572
$ws.Range("I85").Value = 2558.8667
$ws.Range("K85").Value = 2558.8667
$ws.Range("M85").Value = -1310.8667
$ws.Range("H105").Value = 122995
$ws.Range("J105").Value = 122995
$ws.Range("L105").Value = 122995
$ws.Range("N105").Value = -129983
$ws.Range("H136").Value = 5069.52
$ws.Range("I136").Value = 5407.35
$ws.Range("K136").Value = 16222.05
$ws.Range("M136").Value = -13672.05
$ws.Range("N139").ClearContents()
$ws.Range("H139").Value = 74998
$ws.Range("I139").Value = 74998
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 74998
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -69858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 559840
$ws.Range("I2").Value = 2651.4285
$ws.Range("J2").Value = 2510000
$ws.Range("K2").Value = 2651.4285
$ws.Range("L2").Value = 2510000
$ws.Range("M2").Value = -2539.4285
$ws.Range("N2").Value = -2510224
$ws.Range("H113").Value = 603.8333
$ws.Range("I113").Value = 434.75
$ws.Range("J113").Value = 942
$ws.Range("K113").Value = 1304.25
$ws.Range("L113").Value = 2826
$ws.Range("M113").Value = 865.75
$ws.Range("N113").Value = -7166
$ws.Range("H132").Value = 2512.7058
$ws.Range("I132").Value = 2342.3333
$ws.Range("J132").Value = 2921.6
$ws.Range("K132").Value = 7026.999899999999
$ws.Range("L132").Value = 8764.799999999999
$ws.Range("M132").Value = -4496.999899999999
$ws.Range("N132").Value = -13824.8
$ws.Range("H135").Value = 77999.5
$ws.Range("J135").Value = 77999.5
$ws.Range("L135").Value = 77999.5
$ws.Range("N135").Value = -88139.5

